$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Lama1"
$ws.Cells.Item(2, 3).Value = "Rpsa"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.5587383333333333
$ws.Cells.Item(2, 8).Value = 1.676215
$ws.Cells.Item(2, 9).Value = 0.8486764927018626
$ws.Cells.Item(2, 10).Value = 0.8937587278261895
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 75.74527233333333
$ws.Cells.Item(2, 14).Value = 227.235817
$ws.Cells.Item(2, 15).Value = 0.08010992451585024
$ws.Cells.Item(2, 16).Value = 0.08641738889881904
$ws.Cells.Item(2, 17).Value = 42.32178722140611
$ws.Cells.Item(2, 18).Value = 380.896084992655
$ws.Cells.Item(2, 19).Value = 0.06798740976872274
$ws.Cells.Item(2, 20).Value = 0.07723629556426957

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Lama1"
$ws.Cells.Item(3, 3).Value = "Rpsa"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.5587383333333333
$ws.Cells.Item(3, 8).Value = 1.676215
$ws.Cells.Item(3, 9).Value = 0.8486764927018626
$ws.Cells.Item(3, 10).Value = 0.8937587278261895
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 316.1112466666667
$ws.Cells.Item(3, 14).Value = 948.33374
$ws.Cells.Item(3, 15).Value = 0.3343264513940333
$ws.Cells.Item(3, 16).Value = 0.3606496840920618
$ws.Cells.Item(3, 17).Value = 176.6234711104555
$ws.Cells.Item(3, 18).Value = 1589.6112399941
$ws.Cells.Item(3, 19).Value = 0.2837350001865479
$ws.Cells.Item(3, 20).Value = 0.3223338028450383

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Lama1"
$ws.Cells.Item(4, 3).Value = "Rpsa"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.5587383333333333
$ws.Cells.Item(4, 8).Value = 1.676215
$ws.Cells.Item(4, 9).Value = 0.8486764927018626
$ws.Cells.Item(4, 10).Value = 0.8937587278261895
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 178.577067
$ws.Cells.Item(4, 14).Value = 535.7312009999999
$ws.Cells.Item(4, 15).Value = 0.1888671717315399
$ws.Cells.Item(4, 16).Value = 0.203737650838945
$ws.Cells.Item(4, 17).Value = 99.77785278713499
$ws.Cells.Item(4, 18).Value = 898.0006750842149
$ws.Cells.Item(4, 19).Value = 0.1602871288916437
$ws.Cells.Item(4, 20).Value = 0.1820923036241118

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lama1"
$ws.Cells.Item(5, 3).Value = "Rpsa"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.5587383333333333
$ws.Cells.Item(5, 8).Value = 1.676215
$ws.Cells.Item(5, 9).Value = 0.8486764927018626
$ws.Cells.Item(5, 10).Value = 0.8937587278261895
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 207.0351715
$ws.Cells.Item(5, 14).Value = 414.070343
$ws.Cells.Item(5, 15).Value = 0.2189651109577207
$ws.Cells.Item(5, 16).Value = 0.1574702365802588
$ws.Cells.Item(5, 17).Value = 115.6784866652908
$ws.Cells.Item(5, 18).Value = 694.070919991745
$ws.Cells.Item(5, 19).Value = 0.1858305423916726
$ws.Cells.Item(5, 20).Value = 0.1407403983164612

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lama1"
$ws.Cells.Item(6, 3).Value = "Rpsa"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.5587383333333333
$ws.Cells.Item(6, 8).Value = 1.676215
$ws.Cells.Item(6, 9).Value = 0.8486764927018626
$ws.Cells.Item(6, 10).Value = 0.8937587278261895
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 168.0479533333333
$ws.Cells.Item(6, 14).Value = 504.14386
$ws.Cells.Item(6, 15).Value = 0.1777313414008557
$ws.Cells.Item(6, 16).Value = 0.1917250395899155
$ws.Cells.Item(6, 17).Value = 93.89483336554444
$ws.Cells.Item(6, 18).Value = 845.0535002899001
$ws.Cells.Item(6, 19).Value = 0.1508364114632756
$ws.Cells.Item(6, 20).Value = 0.1713559274763087

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Lama1"
$ws.Cells.Item(7, 3).Value = "Rpsa"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.09962599999999999
$ws.Cells.Item(7, 8).Value = 0.199252
$ws.Cells.Item(7, 9).Value = 0.1513235072981373
$ws.Cells.Item(7, 10).Value = 0.1062412721738106
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 75.74527233333333
$ws.Cells.Item(7, 14).Value = 227.235817
$ws.Cells.Item(7, 15).Value = 0.08010992451585024
$ws.Cells.Item(7, 16).Value = 0.08641738889881904
$ws.Cells.Item(7, 17).Value = 7.546198501480666
$ws.Cells.Item(7, 18).Value = 45.277191008884
$ws.Cells.Item(7, 19).Value = 0.01212251474712749
$ws.Cells.Item(7, 20).Value = 0.009181093334549471

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Lama1"
$ws.Cells.Item(8, 3).Value = "Rpsa"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.09962599999999999
$ws.Cells.Item(8, 8).Value = 0.199252
$ws.Cells.Item(8, 9).Value = 0.1513235072981373
$ws.Cells.Item(8, 10).Value = 0.1062412721738106
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 316.1112466666667
$ws.Cells.Item(8, 14).Value = 948.33374
$ws.Cells.Item(8, 15).Value = 0.3343264513940333
$ws.Cells.Item(8, 16).Value = 0.3606496840920618
$ws.Cells.Item(8, 17).Value = 31.49289906041333
$ws.Cells.Item(8, 18).Value = 188.95739436248
$ws.Cells.Item(8, 19).Value = 0.05059145120748535
$ws.Cells.Item(8, 20).Value = 0.03831588124702354

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Lama1"
$ws.Cells.Item(9, 3).Value = "Rpsa"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.09962599999999999
$ws.Cells.Item(9, 8).Value = 0.199252
$ws.Cells.Item(9, 9).Value = 0.1513235072981373
$ws.Cells.Item(9, 10).Value = 0.1062412721738106
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 178.577067
$ws.Cells.Item(9, 14).Value = 535.7312009999999
$ws.Cells.Item(9, 15).Value = 0.1888671717315399
$ws.Cells.Item(9, 16).Value = 0.203737650838945
$ws.Cells.Item(9, 17).Value = 17.79091887694199
$ws.Cells.Item(9, 18).Value = 106.745513261652
$ws.Cells.Item(9, 19).Value = 0.02858004283989624
$ws.Cells.Item(9, 20).Value = 0.02164534721483314

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Lama1"
$ws.Cells.Item(10, 3).Value = "Rpsa"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.09962599999999999
$ws.Cells.Item(10, 8).Value = 0.199252
$ws.Cells.Item(10, 9).Value = 0.1513235072981373
$ws.Cells.Item(10, 10).Value = 0.1062412721738106
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 207.0351715
$ws.Cells.Item(10, 14).Value = 414.070343
$ws.Cells.Item(10, 15).Value = 0.2189651109577207
$ws.Cells.Item(10, 16).Value = 0.1574702365802588
$ws.Cells.Item(10, 17).Value = 20.626085995859
$ws.Cells.Item(10, 18).Value = 82.50434398343599
$ws.Cells.Item(10, 19).Value = 0.03313456856604811
$ws.Cells.Item(10, 20).Value = 0.01672983826379762

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Lama1"
$ws.Cells.Item(11, 3).Value = "Rpsa"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.09962599999999999
$ws.Cells.Item(11, 8).Value = 0.199252
$ws.Cells.Item(11, 9).Value = 0.1513235072981373
$ws.Cells.Item(11, 10).Value = 0.1062412721738106
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 168.0479533333333
$ws.Cells.Item(11, 14).Value = 504.14386
$ws.Cells.Item(11, 15).Value = 0.1777313414008557
$ws.Cells.Item(11, 16).Value = 0.1917250395899155
$ws.Cells.Item(11, 17).Value = 16.74194539878667
$ws.Cells.Item(11, 18).Value = 100.45167239272
$ws.Cells.Item(11, 19).Value = 0.02689492993758013
$ws.Cells.Item(11, 20).Value = 0.02036911211360682
